$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.254.90"
$ws.Range("E2").Value2 = "  -4.42%  "
$ws.Range("D3").Value2 = "3.341.95"
$ws.Range("E3").Value2 = "  -5.39%  "
$ws.Range("E4").Value2 = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "561.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "182.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -6.43%  "
$ws.Range("E7").Value2 = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -2.82%  "
$ws.Range("D9").Value2 = "3.334.92"
$ws.Range("E9").Value2 = "  -5.27%  "
$ws.Range("E10").Value2 = "  -8.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.587"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -6.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "47.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -8.12%  "
$ws.Range("E13").Value2 = "  -6.66%  "
$ws.Range("D14").Value2 = "3.874.66"
$ws.Range("E14").Value2 = "  -5.37%  "
$ws.Range("E15").Value2 = "  -6.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "603.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -9.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "18.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -1.54%  "
$ws.Range("D18").Value2 = "66.286.88"
$ws.Range("E18").Value2 = "  -4.49%  "
$ws.Range("D19").Value2 = "3.338.47"
$ws.Range("E19").Value2 = "  -5.64%  "
$ws.Range("E20").Value2 = "  -3.86%  "
$ws.Range("E21").Value2 = "  -8.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.906"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -5.88%  "
$ws.Range("E23").Value2 = "  -7.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "5.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -4.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "100.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "4.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -7.27%  "
$ws.Range("E27").Value2 = "  -0.01%  "
$ws.Range("E28").Value2 = "  -7.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "9.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -8.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "8.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -8.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "30.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -7.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "6.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -7.15%  "
$ws.Range("E33").Value2 = "  -14.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "11.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -5.93%  "
$ws.Range("E35").Value2 = "  -5.37%  "
$ws.Range("D36").Value2 = "3.799.91"
$ws.Range("E36").Value2 = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "535.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +6.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "57.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -5.31%  "
$ws.Range("D41").Value2 = "0.0₃0714"
$ws.Range("E41").Value2 = "  -11.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "2.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -7.89%  "
$ws.Range("E43").Value2 = "  -6.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.341"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -7.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "31.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -7.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0414"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -7.43%  "
$ws.Range("E47").Value2 = "  -4.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "3.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +14.57%  "
$ws.Range("B49").Value2 = "Stellar"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.129"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -4.80%  "
$ws.Range("B50").Value2 = "ThetaToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -8.28%  "
$ws.Range("E51").Value2 = "  -0.24%  "
